$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove F6 and G6 entirely (ag-test columns not yet populated for this date)
$ws.Range("F6:G6").ClearContents()

# Update AgTests (F) / AgPosit (G) values for revised rows
$ws.Range("F229").Value = 881
$ws.Range("G229").Value = 70
$ws.Range("F271").Value = 45775
$ws.Range("G271").Value = 1732
$ws.Range("F272").Value = 30422
$ws.Range("G272").Value = 1631
$ws.Range("F273").Value = 31409
$ws.Range("G273").Value = 1642
$ws.Range("F274").Value = 28118
$ws.Range("G274").Value = 1277
$ws.Range("F275").Value = 30354
$ws.Range("G275").Value = 1275
$ws.Range("F276").Value = 11428
$ws.Range("F277").Value = 3379
$ws.Range("G277").Value = 130
$ws.Range("F278").Value = 30518
$ws.Range("G278").Value = 2108
$ws.Range("F279").Value = 42730
$ws.Range("G279").Value = 3027
$ws.Range("F280").Value = 34800
$ws.Range("G280").Value = 2318
$ws.Range("F281").Value = 46020
$ws.Range("G281").Value = 3154
$ws.Range("F288").Value = 59282
$ws.Range("F301").Value = 72203
$ws.Range("F302").Value = 78622
$ws.Range("F306").Value = 75382
$ws.Range("G306").Value = 7678
$ws.Range("F307").Value = 75431
$ws.Range("G307").Value = 6326
$ws.Range("F310").Value = 79231
$ws.Range("F311").Value = 61502
$ws.Range("G311").Value = 1926
$ws.Range("F312").Value = 28133
$ws.Range("F313").Value = 75764
$ws.Range("F314").Value = 64292
$ws.Range("F315").Value = 56373
$ws.Range("F318").Value = 49339
$ws.Range("G318").Value = 1134
$ws.Range("F320").Value = 71645
$ws.Range("G320").Value = 3304
$ws.Range("F321").Value = 89635
$ws.Range("G321").Value = 2661
$ws.Range("F322").Value = 109350
$ws.Range("G322").Value = 2331
$ws.Range("F323").Value = 217076
$ws.Range("F324").Value = 249808
$ws.Range("G324").Value = 2856
$ws.Range("F325").Value = 774495
$ws.Range("G325").Value = 6512
$ws.Range("F326").Value = 418334
$ws.Range("G326").Value = 3819
$ws.Range("F327").Value = 223251
$ws.Range("G327").Value = 2717
$ws.Range("F328").Value = 180844
$ws.Range("F329").Value = 73337
$ws.Range("G329").Value = 1729
$ws.Range("F330").Value = 71381
$ws.Range("G330").Value = 2075
$ws.Range("F331").Value = 153763
$ws.Range("G331").Value = 2706
$ws.Range("F332").Value = 485395
$ws.Range("G332").Value = 4801
$ws.Range("F333").Value = 254953
$ws.Range("G333").Value = 2841
$ws.Range("F334").Value = 192639
$ws.Range("G334").Value = 3495
$ws.Range("F335").Value = 150764
$ws.Range("G335").Value = 3786
$ws.Range("F336").Value = 81552
$ws.Range("G336").Value = 2560
$ws.Range("F337").Value = 103428
$ws.Range("G337").Value = 2890
$ws.Range("F338").Value = 221531
$ws.Range("G338").Value = 3071
$ws.Range("F339").Value = 662672
$ws.Range("G339").Value = 5496
$ws.Range("F340").Value = 382054
$ws.Range("G340").Value = 3281
$ws.Range("F341").Value = 283529
$ws.Range("G341").Value = 3611
$ws.Range("F342").Value = 178727
$ws.Range("F343").Value = 133366
$ws.Range("G343").Value = 2979
$ws.Range("F344").Value = 135539
$ws.Range("G344").Value = 2488
$ws.Range("F345").Value = 292191
$ws.Range("G345").Value = 3323
$ws.Range("F346").Value = 674708
$ws.Range("G346").Value = 4823
$ws.Range("F347").Value = 341965
$ws.Range("G347").Value = 2904
$ws.Range("F348").Value = 232760
$ws.Range("G348").Value = 3245
$ws.Range("F349").Value = 159418
$ws.Range("G349").Value = 2758
$ws.Range("F350").Value = 127016
$ws.Range("G350").Value = 2785
$ws.Range("F351").Value = 150635
$ws.Range("G351").Value = 2820
$ws.Range("F352").Value = 307430
$ws.Range("G352").Value = 3546
$ws.Range("F353").Value = 723814
$ws.Range("G353").Value = 5293
$ws.Range("F355").Value = 221975
$ws.Range("G355").Value = 3437
$ws.Range("F356").Value = 160043
$ws.Range("G356").Value = 2883
$ws.Range("F357").Value = 138485
$ws.Range("G357").Value = 3019
$ws.Range("F358").Value = 158776
$ws.Range("G358").Value = 2608
$ws.Range("F359").Value = 321312
$ws.Range("G359").Value = 3339
$ws.Range("F360").Value = 749473
$ws.Range("G360").Value = 5137
$ws.Range("F368").Value = 341507
$ws.Range("G368").Value = 2273
$ws.Range("F369").Value = 234938
$ws.Range("F375").Value = 346616
$ws.Range("G375").Value = 1847
$ws.Range("F376").Value = 221074
$ws.Range("G376").Value = 2216
$ws.Range("F380").Value = 344838
$ws.Range("F381").Value = 745787
$ws.Range("G381").Value = 2690
$ws.Range("F383").Value = 220776
$ws.Range("F384").Value = 172037
$ws.Range("F386").Value = 182753
$ws.Range("G386").Value = 1361
$ws.Range("F391").Value = 176972
$ws.Range("F393").Value = 307442
$ws.Range("G393").Value = 1232
$ws.Range("F394").Value = 166122
$ws.Range("F395").Value = 750082
$ws.Range("G395").Value = 1956
$ws.Range("F398").Value = 298438
$ws.Range("G398").Value = 1471
$ws.Range("F400").Value = 150263
$ws.Range("G400").Value = 758
$ws.Range("F401").Value = 273312
$ws.Range("F402").Value = 716593
$ws.Range("G402").Value = 1383
$ws.Range("F403").Value = 350540
$ws.Range("F404").Value = 224583
$ws.Range("F405").Value = 173767
$ws.Range("G405").Value = 694
$ws.Range("F406").Value = 170555
$ws.Range("F407").Value = 157300
$ws.Range("F408").Value = 301332
$ws.Range("F409").Value = 692081
$ws.Range("G409").Value = 990
$ws.Range("F410").Value = 346053
$ws.Range("F411").Value = 223075
$ws.Range("G411").Value = 819
$ws.Range("F412").Value = 173858
$ws.Range("G412").Value = 639
